$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.513.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.848.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.032'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +2.76%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.027'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4380'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3769'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07395'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8735'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.45'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.855.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.520'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.685'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07165'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.033'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009014'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.027'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.523.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.249'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.075.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.916'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.963'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.250'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09039'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.192'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.30%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7601'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.498'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.78%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.872'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.028'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.147'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01968'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05283'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5139'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.793'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1671'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.707'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.456'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '108.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.703'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4641'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06383'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.46%  '

$ws.Range("E50").Value = '  +2.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.76%  '
